# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Numeric-looking text values (e.g. "218.59") are written with a leading
# apostrophe + ClearFormats() so Excel keeps them as text (matching the
# original inline-string cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '26.106.23'
$ws.Range('E2').Formula = '  -0.55%  '
$ws.Range('D3').Formula = '1.656.68'
$ws.Range('E3').Formula = '  -0.33%  '
$ws.Range('E4').Formula = '  -0.24%  '
$r = $ws.Range('D5')
$r.Formula = '''218.59'
$r.ClearFormats()
$ws.Range('E5').Formula = '  -0.29%  '
$r = $ws.Range('D6')
$r.Formula = '''0.5302'
$r.ClearFormats()
$ws.Range('E6').Formula = '  +1.49%  '
$ws.Range('E7').Formula = '  -0.26%  '
$ws.Range('E8').Formula = '  -1.90%  '
$r = $ws.Range('D9')
$r.Formula = '''0.06341'
$r.ClearFormats()
$ws.Range('E9').Formula = '  -0.07%  '
$r = $ws.Range('D10')
$r.Formula = '''20.46'
$r.ClearFormats()
$ws.Range('E10').Formula = '  -2.97%  '
$r = $ws.Range('D11')
$r.Formula = '''0.07769'
$r.ClearFormats()
$ws.Range('E11').Formula = '  +0.48%  '
$ws.Range('E12').Formula = '  +1.36%  '
$ws.Range('D13').Formula = '1.658.24'
$ws.Range('E13').Formula = '  -0.22%  '
$ws.Range('E14').Formula = '  +0.07%  '
$ws.Range('D15').Formula = '0.0₅8153'
$ws.Range('E15').Formula = '  -0.85%  '
$r = $ws.Range('D16')
$r.Formula = '''65.20'
$r.ClearFormats()
$ws.Range('E16').Formula = '  +0.26%  '
$ws.Range('D17').Formula = '26.104.55'
$ws.Range('E17').Formula = '  -0.62%  '
$ws.Range('E18').Formula = '  -0.30%  '
$r = $ws.Range('D19')
$r.Formula = '''4.555'
$r.ClearFormats()
$ws.Range('E19').Formula = '  -2.22%  '
$r = $ws.Range('D20')
$r.Formula = '''193.25'
$r.ClearFormats()
$ws.Range('E20').Formula = '  -1.30%  '
$r = $ws.Range('D21')
$r.Formula = '''10.08'
$r.ClearFormats()
$ws.Range('E21').Formula = '  -0.66%  '
$r = $ws.Range('D22')
$r.Formula = '''6.033'
$r.ClearFormats()
$ws.Range('E22').Formula = '  -1.00%  '
$r = $ws.Range('D23')
$r.Formula = '''1.004'
$r.ClearFormats()
$ws.Range('E23').Formula = '  -0.39%  '
$r = $ws.Range('D24')
$r.Formula = '''140.01'
$r.ClearFormats()
$ws.Range('E24').Formula = '  +0.52%  '
$ws.Range('E25').Formula = '  -0.06%  '
$ws.Range('E26').Formula = '  +0.45%  '
$r = $ws.Range('D27')
$r.Formula = '''16.18'
$r.ClearFormats()
$ws.Range('E28').Formula = '  +1.16%  '
$r = $ws.Range('D29')
$r.Formula = '''0.05947'
$r.ClearFormats()
$ws.Range('E29').Formula = '  -0.40%  '
$r = $ws.Range('D30')
$r.Formula = '''1.277'
$r.ClearFormats()
$ws.Range('E30').Formula = '  -0.58%  '
$r = $ws.Range('D31')
$r.Formula = '''3.511'
$r.ClearFormats()
$r = $ws.Range('D32')
$r.Formula = '''3.240'
$r.ClearFormats()
$ws.Range('E32').Formula = '  -1.67%  '
$r = $ws.Range('D33')
$r.Formula = '''1.557'
$r.ClearFormats()
$ws.Range('E33').Formula = '  -4.66%  '
$r = $ws.Range('D34')
$r.Formula = '''0.9506'
$r.ClearFormats()
$ws.Range('E34').Formula = '  -3.35%  '
$r = $ws.Range('D35')
$r.Formula = '''2.413'
$r.ClearFormats()
$ws.Range('E35').Formula = '  -0.47%  '
$ws.Range('E36').Formula = '  -0.33%  '
$r = $ws.Range('D37')
$r.Formula = '''0.5657'
$r.ClearFormats()
$ws.Range('E37').Formula = '  -4.09%  '
$ws.Range('E38').Formula = '  +0.77%  '
$r = $ws.Range('D39')
$r.Formula = '''5.836'
$r.ClearFormats()
$ws.Range('E39').Formula = '  -3.09%  '
$r = $ws.Range('D40')
$r.Formula = '''0.8453'
$r.ClearFormats()
$ws.Range('E40').Formula = '  -1.74%  '
$ws.Range('E41').Formula = '  -0.19%  '
$r = $ws.Range('D42')
$r.Formula = '''101.49'
$r.ClearFormats()
$ws.Range('D43').Formula = '1.013.07'
$ws.Range('E43').Formula = '  -1.64%  '
$ws.Range('D44').Formula = '1.800.89'
$r = $ws.Range('D45')
$r.Formula = '''57.16'
$r.ClearFormats()
$ws.Range('E45').Formula = '  -0.50%  '
$ws.Range('E46').Formula = '  -5.02%  '
$ws.Range('E47').Formula = '  -1.31%  '
$ws.Range('E48').Formula = '  +1.40%  '
$ws.Range('B49').Formula = 'RenderToken'
$ws.Range('C49').Formula = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range('D49')
$r.Formula = '''1.472'
$r.ClearFormats()
$ws.Range('E49').Formula = '  +0.34%  '
$ws.Range('B50').Formula = 'Cronos'
$ws.Range('C50').Formula = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r = $ws.Range('D50')
$r.Formula = '''0.05154'
$r.ClearFormats()
$ws.Range('E50').Formula = '  -0.53%  '
$r = $ws.Range('D51')
$r.Formula = '''7.778'
$r.ClearFormats()
$ws.Range('E51').Formula = '  -3.78%  '
